$wb = $excel.ActiveWorkbook
$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

$ws_ALC.Range("H17").Value = 2354.0728
$ws_ALC.Range("J17").Value = 2379.1482
$ws_ALC.Range("L17").Value = 7137.444600000001
$ws_ALC.Range("N17").Value = -7473.444600000001
$ws_ALC.Range("N74").ClearContents()
$ws_ALC.Range("H74").Value = 4900
$ws_ALC.Range("I74").Value = 4900
$ws_ALC.Range("J74").Value = 0
$ws_ALC.Range("K74").Value = 4900
$ws_ALC.Range("L74").Value = 0
$ws_ALC.Range("M74").Value = -3964
$ws_ALC.Range("N77").ClearContents()
$ws_ALC.Range("H77").Value = 4900
$ws_ALC.Range("I77").Value = 4900
$ws_ALC.Range("J77").Value = 0
$ws_ALC.Range("K77").Value = 24500
$ws_ALC.Range("L77").Value = 0
$ws_ALC.Range("M77").Value = -19820
$ws_ALC.Range("H137").Value = 2044.8182
$ws_ALC.Range("I137").Value = 1645.52
$ws_ALC.Range("J137").Value = 2570.2104
$ws_ALC.Range("K137").Value = 4936.559999999999
$ws_ALC.Range("L137").Value = 7710.6312
$ws_ALC.Range("M137").Value = -2386.559999999999
$ws_ALC.Range("N137").Value = -12810.6312
$ws_ARM.Range("H32").Value = 20952.92
$ws_ARM.Range("I32").Value = 22423.037
$ws_ARM.Range("J32").Value = 13161.3
$ws_ARM.Range("K32").Value = 22423.037
$ws_ARM.Range("L32").Value = 13161.3
$ws_ARM.Range("M32").Value = -22136.037
$ws_ARM.Range("N32").Value = -13735.3
$ws_ARM.Range("H61").Value = 7689.512
$ws_ARM.Range("I61").Value = 4324.657
$ws_ARM.Range("J61").Value = 27317.834
$ws_ARM.Range("K61").Value = 4324.657
$ws_ARM.Range("L61").Value = 27317.834
$ws_ARM.Range("M61").Value = -4112.657
$ws_ARM.Range("N61").Value = -27741.834
$ws_ARM.Range("H74").Value = 6119.269
$ws_ARM.Range("I74").Value = 2395.087
$ws_ARM.Range("K74").Value = 2395.087
$ws_ARM.Range("M74").Value = -1521.087
$ws_ARM.Range("H77").Value = 6119.269
$ws_ARM.Range("I77").Value = 2395.087
$ws_ARM.Range("K77").Value = 11975.435
$ws_ARM.Range("M77").Value = -7607.434999999999
$ws_ARM.Range("H122").Value = 1894.8889
$ws_ARM.Range("I122").Value = 1485
$ws_ARM.Range("J122").Value = 2407.25
$ws_ARM.Range("K122").Value = 4455
$ws_ARM.Range("L122").Value = 7221.75
$ws_ARM.Range("M122").Value = -2005
$ws_ARM.Range("N122").Value = -12121.75
$ws_ARM.Range("H132").Value = 13944.909
$ws_ARM.Range("I132").Value = 5708
$ws_ARM.Range("J132").Value = 17033.75
$ws_ARM.Range("K132").Value = 17124
$ws_ARM.Range("L132").Value = 51101.25
$ws_ARM.Range("M132").Value = -14594
$ws_ARM.Range("N132").Value = -56161.25
$ws_ARM.Range("H136").Value = 7689.512
$ws_ARM.Range("I136").Value = 4324.657
$ws_ARM.Range("J136").Value = 27317.834
$ws_ARM.Range("K136").Value = 12973.971
$ws_ARM.Range("L136").Value = 81953.50199999999
$ws_ARM.Range("M136").Value = -10423.971
$ws_ARM.Range("N136").Value = -87053.50199999999
$ws_BSM.Range("H86").Value = 1796.7587
$ws_BSM.Range("I86").Value = 1717.8214
$ws_BSM.Range("J86").Value = 4007
$ws_BSM.Range("K86").Value = 1717.8214
$ws_BSM.Range("L86").Value = 4007
$ws_BSM.Range("M86").Value = -594.8214
$ws_BSM.Range("N86").Value = -6253
$ws_BSM.Range("H89").Value = 1796.7587
$ws_BSM.Range("I89").Value = 1717.8214
$ws_BSM.Range("J89").Value = 4007
$ws_BSM.Range("K89").Value = 8589.107
$ws_BSM.Range("L89").Value = 20035
$ws_BSM.Range("M89").Value = -2973.107
$ws_BSM.Range("N89").Value = -31267
$ws_BSM.Range("H99").Value = 1511.5
$ws_BSM.Range("I99").Value = 1318.091
$ws_BSM.Range("J99").Value = 1815.4286
$ws_BSM.Range("K99").Value = 1318.091
$ws_BSM.Range("L99").Value = 1815.4286
$ws_BSM.Range("M99").Value = 179.9090000000001
$ws_BSM.Range("N99").Value = -4811.4286
$ws_BSM.Range("H107").Value = 1675
$ws_BSM.Range("I107").Value = 1700
$ws_BSM.Range("J107").Value = 1666.6666
$ws_BSM.Range("K107").Value = 1700
$ws_BSM.Range("L107").Value = 1666.6666
$ws_BSM.Range("M107").Value = 220
$ws_BSM.Range("N107").Value = -5506.6666
$ws_CRP.Range("N28").ClearContents()
$ws_CRP.Range("H28").Value = 0
$ws_CRP.Range("J28").Value = 0
$ws_CRP.Range("L28").Value = 0
$ws_CRP.Range("M31").ClearContents()
$ws_CRP.Range("H31").Value = 5080.6665
$ws_CRP.Range("I31").Value = 0
$ws_CRP.Range("J31").Value = 5080.6665
$ws_CRP.Range("K31").Value = 0
$ws_CRP.Range("L31").Value = 5080.6665
$ws_CRP.Range("N31").Value = -5670.6665
$ws_CRP.Range("M34").ClearContents()
$ws_CRP.Range("H34").Value = 5080.6665
$ws_CRP.Range("I34").Value = 0
$ws_CRP.Range("J34").Value = 5080.6665
$ws_CRP.Range("K34").Value = 0
$ws_CRP.Range("L34").Value = 5080.6665
$ws_CRP.Range("N34").Value = -5484.6665
$ws_CRP.Range("H58").Value = 1110153.1
$ws_CRP.Range("I58").Value = 1516179
$ws_CRP.Range("J58").Value = 2809.682
$ws_CRP.Range("K58").Value = 1516179
$ws_CRP.Range("L58").Value = 2809.682
$ws_CRP.Range("M58").Value = -1515976
$ws_CRP.Range("N58").Value = -3215.682
$ws_CRP.Range("H94").Value = 1285.2
$ws_CRP.Range("I94").Value = 936.6667
$ws_CRP.Range("J94").Value = 1434.5714
$ws_CRP.Range("K94").Value = 936.6667
$ws_CRP.Range("L94").Value = 1434.5714
$ws_CRP.Range("M94").Value = -485.6667
$ws_CRP.Range("N94").Value = -2336.5714
$ws_CRP.Range("H136").Value = 1110153.1
$ws_CRP.Range("I136").Value = 1516179
$ws_CRP.Range("J136").Value = 2809.682
$ws_CRP.Range("K136").Value = 4548537
$ws_CRP.Range("L136").Value = 8429.045999999998
$ws_CRP.Range("M136").Value = -4545987
$ws_CRP.Range("N136").Value = -13529.046
$ws_CUL.Range("H6").Value = 74.8
$ws_CUL.Range("I6").Value = 73.666664
$ws_CUL.Range("J6").Value = 76.5
$ws_CUL.Range("K6").Value = 220.999992
$ws_CUL.Range("L6").Value = 229.5
$ws_CUL.Range("M6").Value = -107.999992
$ws_CUL.Range("N6").Value = -455.5
$ws_CUL.Range("H18").Value = 393.22223
$ws_CUL.Range("I18").Value = 328.11765
$ws_CUL.Range("K18").Value = 984.3529500000001
$ws_CUL.Range("M18").Value = -815.3529500000001
$ws_CUL.Range("H129").Value = 2380.3333
$ws_CUL.Range("I129").Value = 3384.2856
$ws_CUL.Range("J129").Value = 1878.3572
$ws_CUL.Range("K129").Value = 10152.8568
$ws_CUL.Range("L129").Value = 5635.071599999999
$ws_CUL.Range("M129").Value = -5152.856800000001
$ws_CUL.Range("N129").Value = -15635.0716
$ws_CUL.Range("H131").Value = 1569.1351
$ws_CUL.Range("I131").Value = 3268
$ws_CUL.Range("J131").Value = 1303.6875
$ws_CUL.Range("K131").Value = 9804
$ws_CUL.Range("L131").Value = 3911.0625
$ws_CUL.Range("M131").Value = -4764
$ws_CUL.Range("N131").Value = -13991.0625
$ws_CUL.Range("H132").Value = 1526
$ws_CUL.Range("I132").Value = 1600.6154
$ws_CUL.Range("J132").Value = 1472.1111
$ws_CUL.Range("K132").Value = 14405.5386
$ws_CUL.Range("L132").Value = 13248.9999
$ws_CUL.Range("M132").Value = -11875.5386
$ws_CUL.Range("N132").Value = -18308.9999
$ws_CUL.Range("H133").Value = 4805.7144
$ws_CUL.Range("I133").Value = 4372.5
$ws_CUL.Range("J133").Value = 4979
$ws_CUL.Range("K133").Value = 13117.5
$ws_CUL.Range("L133").Value = 14937
$ws_CUL.Range("M133").Value = -8057.5
$ws_CUL.Range("N133").Value = -25057
$ws_CUL.Range("H134").Value = 3860.4546
$ws_CUL.Range("I134").Value = 3358.0908
$ws_CUL.Range("K134").Value = 10074.2724
$ws_CUL.Range("M134").Value = -5004.2724
$ws_CUL.Range("H137").Value = 17980.445
$ws_CUL.Range("I137").Value = 1665
$ws_CUL.Range("J137").Value = 22642
$ws_CUL.Range("K137").Value = 4995
$ws_CUL.Range("L137").Value = 67926
$ws_CUL.Range("M137").Value = 105
$ws_CUL.Range("N137").Value = -78126
$ws_CUL.Range("H138").Value = 12447.546
$ws_CUL.Range("I138").Value = 21198
$ws_CUL.Range("J138").Value = 5155.5
$ws_CUL.Range("K138").Value = 63594
$ws_CUL.Range("L138").Value = 15466.5
$ws_CUL.Range("M138").Value = -58454
$ws_CUL.Range("N138").Value = -25746.5
$ws_CUL.Range("H139").Value = 1356571.1
$ws_CUL.Range("I139").Value = 2430487.8
$ws_CUL.Range("K139").Value = 7291463.399999999
$ws_CUL.Range("M139").Value = -7286323.399999999
$ws_CUL.Range("H140").Value = 2560.0393
$ws_CUL.Range("I140").Value = 2004.2667
$ws_CUL.Range("J140").Value = 3354
$ws_CUL.Range("K140").Value = 6012.800099999999
$ws_CUL.Range("L140").Value = 10062
$ws_CUL.Range("M140").Value = -832.8000999999995
$ws_CUL.Range("N140").Value = -20422
$ws_GSM.Range("H70").Value = 6222.154
$ws_GSM.Range("I70").Value = 5871.7896
$ws_GSM.Range("K70").Value = 5871.7896
$ws_GSM.Range("M70").Value = -5601.7896
$ws_GSM.Range("H73").Value = 6222.154
$ws_GSM.Range("I73").Value = 5871.7896
$ws_GSM.Range("K73").Value = 5871.7896
$ws_GSM.Range("M73").Value = -4935.7896
$ws_GSM.Range("H122").Value = 5071.5
$ws_GSM.Range("I122").Value = 8185.7144
$ws_GSM.Range("J122").Value = 2649.3333
$ws_GSM.Range("K122").Value = 24557.1432
$ws_GSM.Range("L122").Value = 7947.999899999999
$ws_GSM.Range("M122").Value = -22107.1432
$ws_GSM.Range("N122").Value = -12847.9999
$ws_GSM.Range("H132").Value = 6478.5654
$ws_GSM.Range("I132").Value = 2320.45
$ws_GSM.Range("J132").Value = 34199.332
$ws_GSM.Range("K132").Value = 6961.349999999999
$ws_GSM.Range("L132").Value = 102597.996
$ws_GSM.Range("M132").Value = -4431.349999999999
$ws_GSM.Range("N132").Value = -107657.996
$ws_LTW.Range("H16").Value = 1304.5834
$ws_LTW.Range("I16").Value = 915.5
$ws_LTW.Range("J16").Value = 2082.75
$ws_LTW.Range("K16").Value = 915.5
$ws_LTW.Range("L16").Value = 2082.75
$ws_LTW.Range("M16").Value = -745.5
$ws_LTW.Range("N16").Value = -2422.75
$ws_LTW.Range("H76").Value = 27716
$ws_LTW.Range("J76").Value = 27716
$ws_LTW.Range("L76").Value = 27716
$ws_LTW.Range("N76").Value = -28392
$ws_LTW.Range("H79").Value = 27716
$ws_LTW.Range("J79").Value = 27716
$ws_LTW.Range("L79").Value = 27716
$ws_LTW.Range("N79").Value = -30056
$ws_LTW.Range("H122").Value = 6906.768
$ws_LTW.Range("I122").Value = 6564.973
$ws_LTW.Range("J122").Value = 7572.3687
$ws_LTW.Range("K122").Value = 19694.919
$ws_LTW.Range("L122").Value = 22717.1061
$ws_LTW.Range("M122").Value = -17244.919
$ws_LTW.Range("N122").Value = -27617.1061
$ws_LTW.Range("H136").Value = 2938.4102
$ws_LTW.Range("I136").Value = 1542.6072
$ws_LTW.Range("K136").Value = 4627.821599999999
$ws_LTW.Range("M136").Value = -2077.821599999999
$ws_WVR.Range("H81").Value = 2211.7
$ws_WVR.Range("J81").Value = 3183.3333
$ws_WVR.Range("L81").Value = 6366.6666
$ws_WVR.Range("N81").Value = -8488.6666
$ws_WVR.Range("H84").Value = 2211.7
$ws_WVR.Range("J84").Value = 3183.3333
$ws_WVR.Range("L84").Value = 31833.333
$ws_WVR.Range("N84").Value = -42441.333
$ws_WVR.Range("H122").Value = 3749.6667
$ws_WVR.Range("I122").Value = 2562.5
$ws_WVR.Range("J122").Value = 4699.4
$ws_WVR.Range("K122").Value = 7687.5
$ws_WVR.Range("L122").Value = 14098.2
$ws_WVR.Range("M122").Value = -5237.5
$ws_WVR.Range("N122").Value = -18998.2
$ws_WVR.Range("H135").Value = 166703040
$ws_WVR.Range("J135").Value = 166703040
$ws_WVR.Range("L135").Value = 166703040
$ws_WVR.Range("N135").Value = -166713180

Write-Host "Applied all changes"
